# Auto-sync batch FINAL (AUTO-TIMEOUT)
# Appends 10 new transaction rows (57-66) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(
    "2026-02-20 11:55:04",
    "2026-02-20 14:59:36",
    "2026-02-20 15:38:46",
    "2026-02-20 13:42:55",
    "2026-02-20 14:58:01",
    "2026-02-20 15:57:46",
    "2026-02-20 11:48:34",
    "2026-02-20 16:09:04",
    "2026-02-20 17:59:48",
    "2026-02-20 13:59:19"
)

$numbers = @(
    "237654101067",
    "237675453374",
    "237680039383",
    "237678973363",
    "237681663743",
    "237651646213",
    "237652194260",
    "237671615641",
    "237653816480",
    "237673593310"
)

$names = @(
    "MAKUETCHE TCHEHGHIE CELINE GIRESSE CHIC MOBILE SARL",
    "ABEL MOUNTAPMBEME",
    "SPECTRUM LTDLA CBOX R0 CEDRICK MARCIALLE WANDJI",
    "MOSSU TAGNE ANNE FLORE TOP MOBIL",
    "LA NEGRESSE SARL FONGA SINTCHA YOLANDE MIREILLE",
    "AMADOU AHIJO ETS MOBILE FINANCIAL SERVICES MFS",
    "CRISTELLE DIANE TCHAHANE",
    "BEGO FOGUE CHRISTELLE KAMILAH CONNECTION GROUP",
    "BERYL NAKOMA TOUFOIN TOP MOBIL TELECOM",
    "JEAN JACQUES YENDJE"
)

$balances = @(82435, 333287, 15241, 217587, 335141, 4643, 436173, 7167, 1929172, 123365)

$startRow = 57
$count = $dates.Count

# Populate column by column so the shared-string table is built in the same
# order the source workbook used (all dates, then all numbers, then all
# names), matching how the rows were authored.
for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $dates[$i]
}

for ($i = 0; $i -lt $count; $i++) {
    # Column B holds long digit-string account numbers that must be stored
    # as text (like the pre-existing rows), not auto-converted to numbers.
    # Force text via NumberFormat, assign, then reset the format back to
    # the default style so no extra formatting sticks to the cell.
    $cellB = $ws.Cells.Item($startRow + $i, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $numbers[$i]
    $cellB.Style = "Normal"
}

for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $names[$i]
}

for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($startRow + $i, 4).Value = $balances[$i]
}
